$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: "Group 1 from Excel" -> "Group UI Automation"
$ws.Range("B2").Value = "Group UI Automation"

# Widen column B (closest the engine's ColumnWidth quantization can reach to 42.140625)
$ws.Columns.Item(2).ColumnWidth = 41.3

# Move the active selection to C2
$ws.Range("C2").Select()
